# Forms the consolidated report: the "Absent" column (H) for each
# attendance row (rows 3-21) is set to the complement of the "Real"
# column (E): Absent = 1 when Real = 0, and Absent = 0 when Real = 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 3; $row -le 21; $row++) {
    $real = $ws.Cells.Item($row, 5).Value()   # column E = "Real"
    if ($real -eq 1) {
        $ws.Cells.Item($row, 8).Value = 0     # column H = "Absent"
    } else {
        $ws.Cells.Item($row, 8).Value = 1
    }
}
